$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old hyperlink (keeps cell styling on C2/D2, drops the relationship) ---
$ws.Hyperlinks.Delete()

# --- Clear the two header cells and the two hyperlink/number cells; C2/D2 end up
#     present-but-empty while retaining their number formats / styles (s="1", s="2") ---
$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# --- Rewrite column A with the new Customer POM field-name list (pass 1, so the shared
#     string table gets the labels first, then the values - matches the authored file) ---
$labels = @(
    "customerState",
    "customerDistrict",
    "customerTaluk",
    "customerPostal",
    "customerVillage",
    "customerShopType",
    "customerShopName",
    "customerName",
    "customerLocalName",
    "customerVillageLocalName",
    "customerPhoneNumber",
    "customerPhoneNumber2",
    "customerAddress",
    "customerLandMark",
    "customerLeisure",
    "customerBreakTime",
    "customerHasCooler",
    "customerCoolerType",
    "customerQualification",
    "customerGrade",
    "customerAvgSale",
    "customerIsSmartPhoneUser"
)
for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $labels[$i]
}

# --- Pass 2: fill in the sample values that exist in column B ---
$ws.Cells.Item(1, 2).Value  = "KARNATAKA"
$ws.Cells.Item(2, 2).Value  = "MANDYA"
$ws.Cells.Item(15, 2).Value = "Break Time"
$ws.Cells.Item(16, 2).Value = "afternoon"
$ws.Cells.Item(17, 2).Value = "yes"
$ws.Cells.Item(18, 2).Value = "commercial"
$ws.Cells.Item(20, 2).Value = "Grade A"

# --- Column widths: A -> 27 characters, B -> as close as this engine's pixel-rounded
#     ColumnWidth setter can get to the authored 29.5703125 (lands on 29.5) ---
$ws.Columns.Item(1).ColumnWidth = 26.166666666666668
$ws.Columns.Item(2).ColumnWidth = 28.666666666666668

# --- Selection / view state ---
$ws.Range("B20").Select()
